$wb = $excel.ActiveWorkbook

# --- 1) Append a new row of data to the "Apple Stock" sheet (row 5) ---
$ws1 = $wb.Worksheets.Item("Apple Stock")
$ws1.Range("A5").Value = "2024-12-31 21:00:00"
$ws1.Range("B5").Value = "2025-01-01 02:30:00"
$ws1.Range("C5").Value = 250.42
$ws1.Range("D5").Value = 20784.86
$ws1.Range("E5").Value = -147.74
$ws1.Range("F5").Value = -0.71
$ws1.Range("G5").Value = 21022.24
$ws1.Range("H5").Value = 20702.69
$ws1.Range("I5").Value = 20952.52
$ws1.Range("J5").Value = 20932.6

# --- 2) Add a new sheet "Gold Prices" right after "Apple Stock" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Gold Prices"

# Header row
$ws2.Range("A1").Value = "Timestamp"
$ws2.Range("B1").Value = "24K Price (INR/g)"
$ws2.Range("C1").Value = "22K Price (INR/g)"
$ws2.Range("D1").Value = "18K Price (INR/g)"

# Match the header formatting used on the "Apple Stock" sheet (bold, bordered,
# centered/top-aligned) by copying the format from its header row.
$ws1.Range("A1:D1").Copy()
$ws2.Range("A1:D1").PasteSpecial(-4122)

# Data rows
$ws2.Range("A2").Value = "2025-01-01 05:17:26"
$ws2.Range("B2").Value = 7230.6965
$ws2.Range("C2").Value = 6628.1384
$ws2.Range("D2").Value = 5423.0224

$ws2.Range("A3").Value = "2025-01-01 06:57:12"
$ws2.Range("B3").Value = 7232.3844
$ws2.Range("C3").Value = 6629.6857
$ws2.Range("D3").Value = 5424.2883

$ws2.Range("A4").Value = "2025-01-01 06:57:12"
$ws2.Range("B4").Value = 7232.3844
$ws2.Range("C4").Value = 6629.6857
$ws2.Range("D4").Value = 5424.2883
